# Disable smart-quote autocorrection so literal straight quotes/apostrophes
# survive any text assignment untouched.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (paragraph 1): an empty run, a bold "Meta description" run, and a
#    plain run with the rest of the sentence.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaRng = $metaPara.Range

$metaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t xml:space="preserve">: Read our review of Collapsed Castle Bonus Buy to discover the pros and cons of this medieval-themed slot machine. Play for free and see if it''s the right game for you!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$metaRng.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2. At the end of the document, drop the duplicated bold title
#    paragraph entirely, and turn the old italic "Read our review..."
#    paragraph into the DALLE image prompt (formatting untouched).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
if ($dupTitlePara.Range.Text.TrimEnd() -ne "Play Collapsed Castle Bonus Buy Free Slot - Pros and Cons") {
    Write-Output ("UNEXPECTED paragraph to delete: [" + $dupTitlePara.Range.Text + "]")
}
$dupTitlePara.Range.Delete()

$count2 = $d.Paragraphs.Count
$promptPara = $d.Paragraphs.Item($count2)
$promptRng = $promptPara.Range
if ($promptRng.Text.TrimEnd() -ne "Read our review of Collapsed Castle Bonus Buy to discover the pros and cons of this medieval-themed slot machine. Play for free and see if it's the right game for you!") {
    Write-Output ("UNEXPECTED paragraph to replace: [" + $promptRng.Text + "]")
}
# Exclude the trailing paragraph mark from the replaced text range.
$textOnlyRng = $d.Range($promptRng.Start, $promptRng.End - 1)
$textOnlyRng.Text = 'Prompt: DALLE, please create a cartoon-style feature image for the game "Collapsed Castle Bonus Buy" that features a happy Maya warrior with glasses. The image should incorporate elements of the medieval fantasy theme, such as a castle in the background or treasure in the foreground. Make sure the image is eye-catching and highlights the adventurous nature of the game.'

Write-Output "Done"
